$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "68.488.55"
$ws.Range("E2").Value = "  -1.17%  "
$ws.Range("D3").Value = "3.832.86"
$ws.Range("E3").Value = "  +2.29%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").Value = "'600.25"
$ws.Range("E5").Value = "  -0.35%  "
$ws.Range("D6").Value = "'163.57"
$ws.Range("E6").Value = "  -2.69%  "
$ws.Range("D7").Value = "3.834.44"
$ws.Range("E7").Value = "  +2.39%  "
$ws.Range("E8").Value = "  +0.12%  "
$ws.Range("D9").Value = "'0.529"
$ws.Range("E9").Value = "  -2.30%  "
$ws.Range("D10").Value = "'0.166"
$ws.Range("E10").Value = "  -2.24%  "
$ws.Range("D11").Value = "'6.34"
$ws.Range("E11").Value = "  -0.67%  "
$ws.Range("E12").Value = "  -0.55%  "
$ws.Range("D13").Value = "'36.79"
$ws.Range("E13").Value = "  -3.92%  "
$ws.Range("E14").Value = "  -1.74%  "
$ws.Range("D15").Value = "4.477.70"
$ws.Range("E15").Value = "  +2.42%  "
$ws.Range("D16").Value = "3.845.34"
$ws.Range("E16").Value = "  +2.75%  "
$ws.Range("D17").Value = "68.694.22"
$ws.Range("E17").Value = "  -0.77%  "
$ws.Range("E18").Value = "  +2.18%  "
$ws.Range("E19").Value = "  -0.34%  "
$ws.Range("D20").Value = "'17.13"
$ws.Range("E20").Value = "  -1.62%  "
$ws.Range("D21").Value = "'11.16"
$ws.Range("E21").Value = "  -0.88%  "
$ws.Range("D22").Value = "'485.93"
$ws.Range("E22").Value = "  -1.44%  "
$ws.Range("D23").Value = "'0.717"
$ws.Range("E23").Value = "  -1.73%  "
$ws.Range("E24").Value = "  +6.66%  "
$ws.Range("D25").Value = "'84.07"
$ws.Range("E25").Value = "  -0.91%  "
$ws.Range("E26").Value = "  -2.62%  "
$ws.Range("D27").Value = "'12.10"
$ws.Range("E27").Value = "  -1.78%  "
$ws.Range("B28").Value = "RenderToken"
$ws.Range("C28").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D28").Value = "'10.00"
$ws.Range("E28").Value = "  -0.83%  "
$ws.Range("B29").Value = "Dai"
$ws.Range("C29").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D29").Value = "'0.998"
$ws.Range("E29").Value = "  -0.22%  "
$ws.Range("E30").Value = "  -0.99%  "
$ws.Range("D31").Value = "'7.84"
$ws.Range("E31").Value = "  -4.15%  "
$ws.Range("D32").Value = "3.986.84"
$ws.Range("E32").Value = "  +2.49%  "
$ws.Range("E33").Value = "  -4.22%  "
$ws.Range("D34").Value = "'31.82"
$ws.Range("E34").Value = "  +0.49%  "
$ws.Range("D35").Value = "3.780.23"
$ws.Range("E35").Value = "  +2.78%  "
$ws.Range("E36").Value = "  -1.40%  "
$ws.Range("E37").Value = "  +1.78%  "
$ws.Range("E38").Value = "  -0.07%  "
$ws.Range("D39").Value = "'5.87"
$ws.Range("E39").Value = "  -1.71%  "
$ws.Range("E40").Value = "  -0.01%  "
$ws.Range("E41").Value = "  -2.93%  "
$ws.Range("D42").Value = "'2.96"
$ws.Range("E42").Value = "  -3.23%  "
$ws.Range("D43").Value = "'428.46"
$ws.Range("E43").Value = "  +0.76%  "
$ws.Range("D44").Value = "'48.48"
$ws.Range("E44").Value = "  -0.88%  "
$ws.Range("E45").Value = "  -0.14%  "
$ws.Range("E46").Value = "  -0.01%  "
$ws.Range("D47").Value = "'8.42"
$ws.Range("E47").Value = "  -0.83%  "
$ws.Range("D48").Value = "2.845.00"
$ws.Range("E48").Value = "  +2.03%  "
$ws.Range("D49").Value = "'142.67"
$ws.Range("E49").Value = "  +0.98%  "
$ws.Range("B50").Value = "EnergySwap"
$ws.Range("C50").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D50").Value = "'26.03"
$ws.Range("E50").Value = "  +13.52%  "
$ws.Range("B51").Value = "VeChain"
$ws.Range("C51").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D51").Value = "'0.0358"
$ws.Range("E51").Value = "  +0.66%  "
